$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 8 new rows after the existing product row (row 7), pushing the
#        totals row (old row 8) and footer row (old row 9) down to rows 16/17 ---
$ws.Range("8:15").Insert()

# Copy the formatting of the template product row (row 7) onto the new rows
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Re-create the merges for the new rows (inserting rows does not copy merges) ---
# First drop + redo the merges that already existed on the (now shifted) totals/footer
# rows so that the merge collection ends up in the same order Excel would produce.
$ws.Range("P16:Q16").UnMerge()
$ws.Range("A17:F17").UnMerge()
$ws.Range("G17:I17").UnMerge()
$ws.Range("K17:Q17").UnMerge()

for ($r = 8; $r -le 15; $r++) {
  $ws.Range("A$r`:B$r").Merge()
  $ws.Range("C$r`:G$r").Merge()
  $ws.Range("H$r`:K$r").Merge()
  $ws.Range("L$r`:M$r").Merge()
  $ws.Range("N$r`:O$r").Merge()
}

$ws.Range("P16:Q16").Merge()
$ws.Range("A17:F17").Merge()
$ws.Range("G17:I17").Merge()
$ws.Range("K17:Q17").Merge()

# --- 3. Fill in the product rows (7 through 15) ---
$products = @(
  @{ Row = 7;  Idx = 1; Name = "B-COM I.M./I.V. 6 AMP";                Stock = "2:2";  Reorder = "1"; Price = "48.00"; SellPrice = "7.6800";  Txns = "0:1" },
  @{ Row = 8;  Idx = 2; Name = "DANSET 4MG/2ML 3 AMP";                 Stock = "1:1";  Reorder = "1"; Price = "82.50"; SellPrice = "27.2250"; Txns = "0:1" },
  @{ Row = 9;  Idx = 3; Name = "DECLOPHEN 75MG/3ML 3 AMPOULES";        Stock = "3:3";  Reorder = "1"; Price = "36.00"; SellPrice = "11.8800"; Txns = "0:1" },
  @{ Row = 10; Idx = 4; Name = "DEXAMETHASONE-AMRIYA 8MG/2ML 3 AMP.";  Stock = "2:0";  Reorder = "1"; Price = "36.00"; SellPrice = "11.8800"; Txns = "0:1" },
  @{ Row = 11; Idx = 5; Name = "جهاز محلول ";                          Stock = "53:0"; Reorder = "0"; Price = "20.00"; SellPrice = "20.0000"; Txns = "1:0" },
  @{ Row = 12; Idx = 6; Name = "سرنجات 3 سم";                          Stock = "0:0";  Reorder = "0"; Price = "2.00";  SellPrice = "2.0000";  Txns = "1:0" },
  @{ Row = 13; Idx = 7; Name = "سرنجات 5 سم";                          Stock = "0:0";  Reorder = "0"; Price = "3.00";  SellPrice = "3.0000";  Txns = "1:0" },
  @{ Row = 14; Idx = 8; Name = "محلول ملح";                            Stock = "6:0";  Reorder = "0"; Price = "24.00"; SellPrice = "24.0000"; Txns = "1:0" },
  @{ Row = 15; Idx = 9; Name = "مناديل سولو سحب صغيره";                Stock = "42:0"; Reorder = "0"; Price = "35.00"; SellPrice = "35.0000"; Txns = "1:0" }
)

foreach ($p in $products) {
  $r = $p.Row
  $ws.Cells.Item($r, 1).Value = $p.Idx                 # A: م (serial number)
  $ws.Cells.Item($r, 3).Value = "'" + $p.Name           # C: الاسم
  $ws.Cells.Item($r, 8).Value = "'" + $p.Stock          # H: الرصيد الحالي
  $ws.Cells.Item($r, 12).Value = "'" + $p.Reorder       # L: حد الطلب
  $ws.Cells.Item($r, 14).Value = "'" + $p.Price         # N: السعر
  $ws.Cells.Item($r, 16).Value = "'" + $p.SellPrice     # P: سعر البيع
  $ws.Cells.Item($r, 17).Value = "'" + $p.Txns          # Q: عدد التعااملات
}

# --- 4. Update the totals row (now row 16) ---
$ws.Cells.Item(16, 16).Value = 142.665

# --- 5. Update the generated-on timestamp in the footer (now row 17) ---
$ws.Cells.Item(17, 1).Value = "Friday, 23 May, 2025 2:46 PM"

Write-Host "Edit complete"
